$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map derived from the source diff (crypto price list refresh).
# All of these columns are stored as text in the workbook (e.g. "3.409.37", "  -1.57%  "),
# so each cell is forced to Text format before the value is written, then restored to
# "General" (matching the original formatting) to avoid Excel auto-converting the strings
# that look like numbers (e.g. "1.00", "0.0000200") into actual numeric values.
$updates = [ordered]@{
    "D2" = "65.534.48"
    "E2" = "  -1.98%  "
    "D3" = "3.400.48"
    "E3" = "  -1.89%  "
    "E4" = "  +0.10%  "
    "D5" = "595.57"
    "E5" = "  -1.49%  "
    "D6" = "141.84"
    "E6" = "  -4.28%  "
    "B7" = "USDC"
    "C7" = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
    "D7" = "1.00"
    "E7" = "  -0.02%  "
    "B8" = "LidoStakedEther"
    "C8" = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
    "D8" = "3.399.74"
    "E8" = "  -1.84%  "
    "D9" = "0.469"
    "E9" = "  -3.01%  "
    "D10" = "0.134"
    "E10" = "  -5.57%  "
    "D11" = "7.88"
    "E11" = "  +5.13%  "
    "D12" = "0.407"
    "E12" = "  -4.20%  "
    "D13" = "3.979.08"
    "E13" = "  -1.79%  "
    "D14" = "0.0000200"
    "E14" = "  -7.13%  "
    "D15" = "29.78"
    "E15" = "  -6.27%  "
    "E16" = "  -0.40%  "
    "D17" = "3.405.15"
    "E17" = "  -1.82%  "
    "D18" = "65.596.35"
    "E18" = "  -1.95%  "
    "D19" = "10.39"
    "E19" = "  +3.86%  "
    "D20" = "6.10"
    "E20" = "  -5.60%  "
    "D21" = "14.78"
    "E21" = "  -3.98%  "
    "D22" = "416.28"
    "E22" = "  -5.42%  "
    "D23" = "0.580"
    "E23" = "  -5.00%  "
    "D24" = "77.46"
    "E24" = "  -1.65%  "
    "E25" = "  +0.07%  "
    "D26" = "3.540.07"
    "E26" = "  -1.72%  "
    "E27" = "  -8.87%  "
    "D28" = "9.27"
    "E28" = "  -5.97%  "
    "D29" = "7.82"
    "E29" = "  -7.41%  "
    "D30" = "2.42"
    "E30" = "  -2.44%  "
    "D31" = "1.00"
    "E31" = "  +0.11%  "
    "D32" = "0.160"
    "E32" = "  -4.76%  "
    "D33" = "1.47"
    "E33" = "  -8.73%  "
    "D34" = "24.48"
    "E34" = "  -3.76%  "
    "D35" = "3.401.71"
    "E35" = "  -1.62%  "
    "E36" = "  -0.06%  "
    "D37" = "1.70"
    "E37" = "  -6.17%  "
    "D38" = "5.55"
    "E38" = "  -8.83%  "
    "D39" = "7.56"
    "E39" = "  -4.64%  "
    "D40" = "1.00"
    "E40" = "  +0.26%  "
    "D41" = "170.21"
    "E41" = "  -2.06%  "
    "D42" = "0.0859"
    "E42" = "  -3.82%  "
    "D43" = "5.06"
    "E43" = "  -6.61%  "
    "D44" = "0.870"
    "E44" = "  -1.61%  "
    "D45" = "1.92"
    "E45" = "  -11.37%  "
    "D46" = "45.46"
    "E46" = "  -1.11%  "
    "D47" = "26.86"
    "E47" = "  -7.90%  "
    "D48" = "1.17"
    "E48" = "  -6.12%  "
    "D49" = "7.07"
    "E49" = "  -5.58%  "
    "D50" = "2.29"
    "E50" = "  -7.33%  "
    "D51" = "0.921"
    "E51" = "  -6.72%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = "General"
}

Write-Output "Updated $($updates.Count) cells"
